$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44176
$ws.Cells.Item(2, 12).Value = 'Primera'
$ws.Cells.Item(2, 13).Value = 300
$ws.Cells.Item(2, 14).Value = 5000
$ws.Cells.Item(2, 15).Value = 6000
$ws.Cells.Item(2, 16).Value = 5500
$ws.Cells.Item(2, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(2, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(2, 19).Value = 3667
$ws.Cells.Item(2, 20).Value = 1.5

$ws.Cells.Item(3, 4).Value = 44553
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 13).Value = 400
$ws.Cells.Item(3, 14).Value = 5000
$ws.Cells.Item(3, 15).Value = 5500
$ws.Cells.Item(3, 16).Value = 5250
$ws.Cells.Item(3, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(3, 18).Value = 'Región del Maule'
$ws.Cells.Item(3, 19).Value = 3500
$ws.Cells.Item(3, 20).Value = 1.5

$ws.Cells.Item(4, 4).Value = 44547
$ws.Cells.Item(4, 12).Value = 'Primera'
$ws.Cells.Item(4, 13).Value = 400
$ws.Cells.Item(4, 14).Value = 5000
$ws.Cells.Item(4, 15).Value = 5500
$ws.Cells.Item(4, 16).Value = 5250
$ws.Cells.Item(4, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(4, 18).Value = 'Región del Maule'
$ws.Cells.Item(4, 19).Value = 3500
$ws.Cells.Item(4, 20).Value = 1.5

$ws.Cells.Item(5, 4).Value = 44551
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 400
$ws.Cells.Item(5, 14).Value = 5000
$ws.Cells.Item(5, 15).Value = 5500
$ws.Cells.Item(5, 16).Value = 5250
$ws.Cells.Item(5, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(5, 18).Value = 'Región del Maule'
$ws.Cells.Item(5, 19).Value = 3500
$ws.Cells.Item(5, 20).Value = 1.5

$ws.Cells.Item(6, 4).Value = 44537
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 400
$ws.Cells.Item(6, 14).Value = 5000
$ws.Cells.Item(6, 15).Value = 5500
$ws.Cells.Item(6, 16).Value = 5250
$ws.Cells.Item(6, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(6, 18).Value = 'Región del Maule'
$ws.Cells.Item(6, 19).Value = 3500
$ws.Cells.Item(6, 20).Value = 1.5

$ws.Cells.Item(7, 4).Value = 44159
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 200
$ws.Cells.Item(7, 14).Value = 6500
$ws.Cells.Item(7, 15).Value = 7000
$ws.Cells.Item(7, 16).Value = 6750
$ws.Cells.Item(7, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(7, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(7, 19).Value = 4500
$ws.Cells.Item(7, 20).Value = 1.5

$ws.Cells.Item(8, 4).Value = 44519
$ws.Cells.Item(8, 12).Value = 'Primera'
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 3700
$ws.Cells.Item(8, 15).Value = 3800
$ws.Cells.Item(8, 16).Value = 3750
$ws.Cells.Item(8, 17).Value = '$/kilo'
$ws.Cells.Item(8, 18).Value = 'Región del Maule'
$ws.Cells.Item(8, 19).Value = 3750
$ws.Cells.Item(8, 20).Value = 1

$ws.Cells.Item(9, 4).Value = 44166
$ws.Cells.Item(9, 12).Value = 'Primera'
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 6000
$ws.Cells.Item(9, 15).Value = 6500
$ws.Cells.Item(9, 16).Value = 6250
$ws.Cells.Item(9, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(9, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(9, 19).Value = 4167
$ws.Cells.Item(9, 20).Value = 1.5

$ws.Cells.Item(10, 4).Value = 44530
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 160
$ws.Cells.Item(10, 14).Value = 3600
$ws.Cells.Item(10, 15).Value = 3700
$ws.Cells.Item(10, 16).Value = 3650
$ws.Cells.Item(10, 17).Value = '$/kilo'
$ws.Cells.Item(10, 18).Value = 'Región del Maule'
$ws.Cells.Item(10, 19).Value = 3650
$ws.Cells.Item(10, 20).Value = 1

$ws.Cells.Item(11, 4).Value = 44516
$ws.Cells.Item(11, 12).Value = 'Primera'
$ws.Cells.Item(11, 13).Value = 80
$ws.Cells.Item(11, 14).Value = 3700
$ws.Cells.Item(11, 15).Value = 3800
$ws.Cells.Item(11, 16).Value = 3750
$ws.Cells.Item(11, 17).Value = '$/kilo'
$ws.Cells.Item(11, 18).Value = 'Región del Maule'
$ws.Cells.Item(11, 19).Value = 3750
$ws.Cells.Item(11, 20).Value = 1

$ws.Cells.Item(12, 4).Value = 44523
$ws.Cells.Item(12, 12).Value = 'Primera'
$ws.Cells.Item(12, 13).Value = 300
$ws.Cells.Item(12, 14).Value = 3700
$ws.Cells.Item(12, 15).Value = 3800
$ws.Cells.Item(12, 16).Value = 3750
$ws.Cells.Item(12, 17).Value = '$/kilo'
$ws.Cells.Item(12, 18).Value = 'Región del Maule'
$ws.Cells.Item(12, 19).Value = 3750
$ws.Cells.Item(12, 20).Value = 1

$ws.Cells.Item(13, 4).Value = 44544
$ws.Cells.Item(13, 12).Value = 'Primera'
$ws.Cells.Item(13, 13).Value = 400
$ws.Cells.Item(13, 14).Value = 5000
$ws.Cells.Item(13, 15).Value = 5500
$ws.Cells.Item(13, 16).Value = 5250
$ws.Cells.Item(13, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(13, 18).Value = 'Región del Maule'
$ws.Cells.Item(13, 19).Value = 3500
$ws.Cells.Item(13, 20).Value = 1.5

$ws.Cells.Item(14, 4).Value = 44169
$ws.Cells.Item(14, 12).Value = 'Primera'
$ws.Cells.Item(14, 13).Value = 400
$ws.Cells.Item(14, 14).Value = 5500
$ws.Cells.Item(14, 15).Value = 6000
$ws.Cells.Item(14, 16).Value = 5750
$ws.Cells.Item(14, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(14, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(14, 19).Value = 3833
$ws.Cells.Item(14, 20).Value = 1.5

$ws.Cells.Item(15, 4).Value = 44162
$ws.Cells.Item(15, 12).Value = 'Primera'
$ws.Cells.Item(15, 13).Value = 100
$ws.Cells.Item(15, 14).Value = 7000
$ws.Cells.Item(15, 15).Value = 7000
$ws.Cells.Item(15, 16).Value = 7000
$ws.Cells.Item(15, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(15, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(15, 19).Value = 4667
$ws.Cells.Item(15, 20).Value = 1.5

$ws.Cells.Item(16, 4).Value = 44162
$ws.Cells.Item(16, 12).Value = 'Segunda'
$ws.Cells.Item(16, 13).Value = 100
$ws.Cells.Item(16, 14).Value = 6500
$ws.Cells.Item(16, 15).Value = 6500
$ws.Cells.Item(16, 16).Value = 6500
$ws.Cells.Item(16, 17).Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Cells.Item(16, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(16, 19).Value = 4333
$ws.Cells.Item(16, 20).Value = 1.5

$ws.Cells.Item(17, 4).Value = 44533
$ws.Cells.Item(17, 12).Value = 'Primera'
$ws.Cells.Item(17, 13).Value = 400
$ws.Cells.Item(17, 14).Value = 3500
$ws.Cells.Item(17, 15).Value = 3600
$ws.Cells.Item(17, 16).Value = 3550
$ws.Cells.Item(17, 17).Value = '$/kilo'
$ws.Cells.Item(17, 18).Value = 'Región del Maule'
$ws.Cells.Item(17, 19).Value = 3550
$ws.Cells.Item(17, 20).Value = 1
